$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 113's "Numero de page" (C113) used to hold "NA"; that value now moves
# down to the newly appended row 114, so C113 becomes blank.
$ws.Range("C113").Value = ""

# Append the new data row (114) that the script produced on its latest run.
# Column A holds a date formatted/stored as plain text (matching every other
# row), so force text formatting before assigning to stop Excel from
# auto-converting the literal "2025-05-20" into a date serial number; then
# restore the default "Normal" style so no stray number-format/style is left
# behind on the cell.
$ws.Range("A114").NumberFormat = "@"
$ws.Range("A114").Value = "2025-05-20"
$ws.Range("A114").Style = "Normal"

$ws.Range("B114").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C114").Value = "NA"
$ws.Range("D114").Value = 1
